# Auto-generated edit script applying the Marilith_Profits.xlsx diff
# Updates H/I/J/K/L/M/N profit-calculation cells across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H76").Value = 4277.6665
$ws.Range("J76").Value = 4750
$ws.Range("L76").Value = 4750
$ws.Range("N76").Value = -5380
$ws.Range("H79").Value = 4277.6665
$ws.Range("J79").Value = 4750
$ws.Range("L79").Value = 4750
$ws.Range("N79").Value = -6934
$ws.Range("H118").Value = 352.875
$ws.Range("I118").Value = 352.875
$ws.Range("K118").Value = 1058.625
$ws.Range("M118").Value = 598.375
$ws.Range("H125").Value = 5134.8
$ws.Range("I125").Value = 3841.6667
$ws.Range("J125").Value = 7074.5
$ws.Range("K125").Value = 34575.0003
$ws.Range("L125").Value = 63670.5
$ws.Range("M125").Value = -32115.0003
$ws.Range("N125").Value = -68590.5
$ws.Range("H129").Value = 1101.4073
$ws.Range("J129").Value = 3411.4
$ws.Range("L129").Value = 10234.2
$ws.Range("N129").Value = -20234.2
$ws.Range("H132").Value = 2921.5334
$ws.Range("I132").Value = 2921.5334
$ws.Range("K132").Value = 8764.600199999999
$ws.Range("M132").Value = -6234.600199999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1593.258
$ws.Range("I32").Value = 1619.7
$ws.Range("K32").Value = 1619.7
$ws.Range("M32").Value = -1332.7
$ws.Range("H61").Value = 1750.8462
$ws.Range("I61").Value = 1772.1666
$ws.Range("K61").Value = 1772.1666
$ws.Range("M61").Value = -1560.1666
$ws.Range("H63").Value = 3475
$ws.Range("H66").Value = 3475
$ws.Range("H88").Value = 3106.923
$ws.Range("I88").Value = 1110
$ws.Range("J88").Value = 3706
$ws.Range("K88").Value = 1110
$ws.Range("L88").Value = 3706
$ws.Range("M88").Value = -704
$ws.Range("N88").Value = -4518
$ws.Range("H91").Value = 3106.923
$ws.Range("I91").Value = 1110
$ws.Range("J91").Value = 3706
$ws.Range("K91").Value = 1110
$ws.Range("L91").Value = 3706
$ws.Range("M91").Value = 294
$ws.Range("N91").Value = -6514
$ws.Range("H132").Value = 1381.5186
$ws.Range("I132").Value = 1306.6086
$ws.Range("K132").Value = 3919.8258
$ws.Range("M132").Value = -1389.8258
$ws.Range("H135").Value = 59000
$ws.Range("J135").Value = 59000
$ws.Range("L135").Value = 59000
$ws.Range("N135").Value = -69140
$ws.Range("H136").Value = 1750.8462
$ws.Range("I136").Value = 1772.1666
$ws.Range("K136").Value = 5316.4998
$ws.Range("M136").Value = -2766.4998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1917.3636
$ws.Range("I86").Value = 1819.1
$ws.Range("K86").Value = 1819.1
$ws.Range("M86").Value = -696.0999999999999
$ws.Range("H89").Value = 1917.3636
$ws.Range("I89").Value = 1819.1
$ws.Range("K89").Value = 9095.5
$ws.Range("M89").Value = -3479.5
$ws.Range("H105").Value = 5054.3
$ws.Range("I105").Value = 5078.143
$ws.Range("J105").Value = 4998.6665
$ws.Range("K105").Value = 5078.143
$ws.Range("L105").Value = 4998.6665
$ws.Range("M105").Value = -3331.143
$ws.Range("N105").Value = -8492.666499999999
$ws.Range("H134").Value = 3652
$ws.Range("I134").Value = 3946.7778
$ws.Range("K134").Value = 11840.3334
$ws.Range("M134").Value = -9305.3334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 371
$ws.Range("I2").Value = 368.33334
$ws.Range("J2").Value = 375
$ws.Range("K2").Value = 368.33334
$ws.Range("L2").Value = 375
$ws.Range("M2").Value = -255.33334
$ws.Range("N2").Value = -601
$ws.Range("H4").Value = 26750
$ws.Range("H31").Value = 2255.5945
$ws.Range("I31").Value = 1712.75
$ws.Range("K31").Value = 1712.75
$ws.Range("M31").Value = -1417.75
$ws.Range("H34").Value = 2255.5945
$ws.Range("I34").Value = 1712.75
$ws.Range("K34").Value = 1712.75
$ws.Range("M34").Value = -1510.75
$ws.Range("H58").Value = 4433.3
$ws.Range("I58").Value = 4216.9443
$ws.Range("K58").Value = 4216.9443
$ws.Range("M58").Value = -4013.9443
$ws.Range("H92").Value = 51629.332
$ws.Range("I92").Value = 36000
$ws.Range("J92").Value = 59444
$ws.Range("K92").Value = 36000
$ws.Range("L92").Value = 59444
$ws.Range("M92").Value = -33504
$ws.Range("N92").Value = -64436
$ws.Range("H132").Value = 2277.8
$ws.Range("I132").Value = 1972.25
$ws.Range("K132").Value = 5916.75
$ws.Range("M132").Value = -3386.75
$ws.Range("H134").Value = 4499.75
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H136").Value = 4433.3
$ws.Range("I136").Value = 4216.9443
$ws.Range("K136").Value = 12650.8329
$ws.Range("M136").Value = -10100.8329

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1789.0869
$ws.Range("I4").Value = 1696.8823
$ws.Range("J4").Value = 2050.3333
$ws.Range("K4").Value = 5090.6469
$ws.Range("L4").Value = 6150.999899999999
$ws.Range("M4").Value = -4978.6469
$ws.Range("N4").Value = -6374.999899999999
$ws.Range("H7").Value = 419.83334
$ws.Range("I7").Value = 150
$ws.Range("J7").Value = 473.8
$ws.Range("K7").Value = 450
$ws.Range("L7").Value = 1421.4
$ws.Range("M7").Value = -338
$ws.Range("N7").Value = -1645.4
$ws.Range("H10").Value = 12.555555
$ws.Range("I10").Value = 12.555555
$ws.Range("K10").Value = 37.666665
$ws.Range("M10").Value = 101.333335
$ws.Range("H13").Value = 185.375
$ws.Range("I13").Value = 352.25
$ws.Range("J13").Value = 18.5
$ws.Range("K13").Value = 1056.75
$ws.Range("L13").Value = 55.5
$ws.Range("M13").Value = -888.75
$ws.Range("N13").Value = -391.5
$ws.Range("H37").Value = 99950
$ws.Range("J37").Value = 99950
$ws.Range("L37").Value = 299850
$ws.Range("N37").Value = -300074
$ws.Range("H56").Value = 11466.75
$ws.Range("I56").Value = 11466.75
$ws.Range("K56").Value = 11466.75
$ws.Range("M56").Value = -10936.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4224.5
$ws.Range("I70").Value = 3499
$ws.Range("J70").Value = 4950
$ws.Range("K70").Value = 3499
$ws.Range("L70").Value = 4950
$ws.Range("M70").Value = -3229
$ws.Range("N70").Value = -5490
$ws.Range("H73").Value = 4224.5
$ws.Range("I73").Value = 3499
$ws.Range("J73").Value = 4950
$ws.Range("K73").Value = 3499
$ws.Range("L73").Value = 4950
$ws.Range("M73").Value = -2563
$ws.Range("N73").Value = -6822
$ws.Range("H98").Value = 8735.625
$ws.Range("J98").Value = 8735.625
$ws.Range("L98").Value = 8735.625
$ws.Range("N98").Value = -14725.625
$ws.Range("H126").Value = 9499.75
$ws.Range("I126").Value = 9499.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 28499.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -26029.25
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 988.7
$ws.Range("I132").Value = 999.2222
$ws.Range("J132").Value = 894
$ws.Range("K132").Value = 2997.6666
$ws.Range("L132").Value = 2682
$ws.Range("M132").Value = -467.6666
$ws.Range("N132").Value = -7742

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1500
$ws.Range("I9").Value = 833.3333
$ws.Range("K9").Value = 833.3333
$ws.Range("M9").Value = -609.3333
$ws.Range("H22").Value = 1038.8334
$ws.Range("I22").Value = 1049.9286
$ws.Range("K22").Value = 1049.9286
$ws.Range("M22").Value = -754.9286
$ws.Range("H27").Value = 1038.8334
$ws.Range("I27").Value = 1049.9286
$ws.Range("K27").Value = 1049.9286
$ws.Range("M27").Value = -942.9286
$ws.Range("H46").Value = 3215
$ws.Range("J46").Value = 4069.2856
$ws.Range("L46").Value = 4069.2856
$ws.Range("N46").Value = -4445.2856
$ws.Range("H55").Value = 258.7647
$ws.Range("I55").Value = 196
$ws.Range("J55").Value = 329.375
$ws.Range("K55").Value = 196
$ws.Range("L55").Value = 329.375
$ws.Range("M55").Value = -23
$ws.Range("N55").Value = -675.375
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 2226.25
$ws.Range("I136").Value = 2226.25
$ws.Range("K136").Value = 6678.75
$ws.Range("M136").Value = -4128.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5000.5
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = 5001
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 5001
$ws.Range("M2").Value = -4888
$ws.Range("N2").Value = -5225
$ws.Range("H132").Value = 1498.75
$ws.Range("I132").Value = 1498.6666
$ws.Range("K132").Value = 4495.9998
$ws.Range("M132").Value = -1965.9998
$ws.Range("H133").Value = 85107.5
$ws.Range("J133").Value = 80143.336
$ws.Range("L133").Value = 80143.336
$ws.Range("N133").Value = -90263.336

